$d = $word.ActiveDocument

# 1) Heading text: "User scenarios:" -> "User scenario:"
$d.Content.Find.Execute("User scenarios:", $false, $false, $false, $false, $false,
                         $true, 1, $false, "User scenario:", 2)

# 2) Collapse paragraphs 2..11 (the blank line + the six "As a user, ..." bullet
#    paragraphs, each separated by a blank paragraph) away, then fill the blank
#    paragraph that is left (old paragraph 12, now paragraph 2 after the
#    deletion) with the new scenario narrative. Old paragraph 13 is left
#    completely untouched as its own blank paragraph right after.
$start = $d.Paragraphs.Item(2).Range.Start
$finish = $d.Paragraphs.Item(11).Range.End
$r = $d.Range($start, $finish)
$r.Delete()

$scenario = "A user will use the spell checker to compare text inputs with a dictionary. When the program starts, the user will see a help screen. The user will be prompted by the program to enter a name of an input text file and the dictionary. The program will then display all the words from the input text file that are not in the dictionary. At this point, the user can decide to either ignore each of these words or to add them to the dictionary. The user will be able to store all text and words before opening another file. In the case that the user inputs the wrong name for a file, an error message will display, and the user can try again. The user can repeat the process of uploading a text file to compare to the dictionary as many times as he/she wishes, and the program will only end when the user chooses to close it."

# Insert the narrative plus a temporary paragraph break in front of the blank
# paragraph's existing (empty) run, then delete that break again. This keeps
# the original empty run intact as its own trailing <w:r/> (matching the
# source formatting) instead of Word silently merging the new text into it.
$p2Start = $d.Paragraphs.Item(2).Range.Start
$insertPoint = $d.Range($p2Start, $p2Start)
$insertPoint.InsertAfter("$scenario`r")

$breakPos = $p2Start + $scenario.Length
$breakRange = $d.Range($breakPos, $breakPos + 1)
$breakRange.Delete()

# 3) "Mock GUI" -> "Mock GUI:"
$d.Content.Find.Execute("Mock GUI", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mock GUI:", 2)
